$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 428.8
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562
$ws.Range("H33").Value = 154
$ws.Range("I33").Value = 145.66667
$ws.Range("K33").Value = 145.66667
$ws.Range("M33").Value = 83.33332999999999
$ws.Range("H38").Value = 304.3
$ws.Range("I38").Value = 227
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 681
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -309
$ws.Range("N38").Value = -3744
$ws.Range("H41").Value = 2000
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1560
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value = 2292.1538
$ws.Range("I58").Value = 254.42857
$ws.Range("J58").Value = 4669.5
$ws.Range("K58").Value = 763.28571
$ws.Range("L58").Value = 14008.5
$ws.Range("M58").Value = -613.28571
$ws.Range("N58").Value = -14308.5
$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496
$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480
$ws.Range("H113").Value = 6245.7144
$ws.Range("I113").Value = 5305.75
$ws.Range("K113").Value = 5305.75
$ws.Range("M113").Value = -2051.75
$ws.Range("H116").Value = 3647.5
$ws.Range("I116").Value = 3647.5
$ws.Range("K116").Value = 3647.5
$ws.Range("M116").Value = -205.5
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 15000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20060
$ws.Range("H141").Value = 2693.6924
$ws.Range("I141").Value = 1859.5555
$ws.Range("J141").Value = 3135.2942
$ws.Range("K141").Value = 5578.666499999999
$ws.Range("L141").Value = 9405.882599999999
$ws.Range("M141").Value = -398.6664999999994
$ws.Range("N141").Value = -19765.8826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1037.1666
$ws.Range("I2").Value = 1164.6
$ws.Range("K2").Value = 1164.6
$ws.Range("M2").Value = -1051.6
$ws.Range("H32").Value = 1908.7742
$ws.Range("I32").Value = 1164.5385
$ws.Range("K32").Value = 1164.5385
$ws.Range("M32").Value = -877.5385000000001
$ws.Range("H101").Value = 65000
$ws.Range("J101").Value = 65000
$ws.Range("L101").Value = 65000
$ws.Range("N101").Value = -71490
$ws.Range("H116").Value = 1037.1666
$ws.Range("I116").Value = 1164.6
$ws.Range("K116").Value = 1164.6
$ws.Range("M116").Value = 1129.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1037.1666
$ws.Range("I3").Value = 1164.6
$ws.Range("K3").Value = 1164.6
$ws.Range("M3").Value = -1050.6
$ws.Range("H94").Value = 1162.2858
$ws.Range("I94").Value = 1042.8889
$ws.Range("K94").Value = 1042.8889
$ws.Range("M94").Value = -591.8888999999999
$ws.Range("H134").Value = 5903.2144
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7855.3335
$ws.Range("J86").Value = 8783
$ws.Range("L86").Value = 8783
$ws.Range("N86").Value = -11029
$ws.Range("H89").Value = 7855.3335
$ws.Range("J89").Value = 8783
$ws.Range("L89").Value = 43915
$ws.Range("N89").Value = -55147
$ws.Range("H99").Value = 3124.625
$ws.Range("I99").Value = 3385.2856
$ws.Range("K99").Value = 3385.2856
$ws.Range("M99").Value = -1887.2856
$ws.Range("H126").Value = 3124.625
$ws.Range("I126").Value = 3385.2856
$ws.Range("K126").Value = 10155.8568
$ws.Range("M126").Value = -7685.856800000001
$ws.Range("H134").Value = 3679.7778
$ws.Range("I134").Value = 2716.5833
$ws.Range("J134").Value = 4450.3335
$ws.Range("K134").Value = 8149.749899999999
$ws.Range("L134").Value = 13351.0005
$ws.Range("M134").Value = -5614.749899999999
$ws.Range("N134").Value = -18421.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 980
$ws.Range("I133").Value = 980
$ws.Range("K133").Value = 2940
$ws.Range("M133").Value = 2120
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H137").Value = 2124.75
$ws.Range("I137").Value = 999.5
$ws.Range("J137").Value = 3250
$ws.Range("K137").Value = 2998.5
$ws.Range("L137").Value = 9750
$ws.Range("M137").Value = 2101.5
$ws.Range("N137").Value = -19950
$ws.Range("H138").Value = 1365.7142
$ws.Range("I138").Value = 1343.3334
$ws.Range("K138").Value = 4030.0002
$ws.Range("M138").Value = 1109.9998
$ws.Range("H139").Value = 2413
$ws.Range("I139").Value = 621.25
$ws.Range("K139").Value = 1863.75
$ws.Range("M139").Value = 3276.25
$ws.Range("H141").Value = 2225.6
$ws.Range("I141").Value = 2225.6
$ws.Range("K141").Value = 6676.799999999999
$ws.Range("M141").Value = -1496.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749.5
$ws.Range("I80").Value = 2999
$ws.Range("K80").Value = 2999
$ws.Range("M80").Value = -2001
$ws.Range("H83").Value = 2749.5
$ws.Range("I83").Value = 2999
$ws.Range("K83").Value = 14995
$ws.Range("M83").Value = -10003
$ws.Range("H104").Value = 31999.5
$ws.Range("J104").Value = 31999.5
$ws.Range("L104").Value = 31999.5
$ws.Range("N104").Value = -38987.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1498.2
$ws.Range("I96").Value = 872.75
$ws.Range("K96").Value = 872.75
$ws.Range("M96").Value = 500.25
$ws.Range("H136").Value = 3647.6
$ws.Range("I136").Value = 3773.7646
$ws.Range("J136").Value = 3379.5
$ws.Range("K136").Value = 11321.2938
$ws.Range("L136").Value = 10138.5
$ws.Range("M136").Value = -8771.2938
$ws.Range("N136").Value = -15238.5

